# Auto-generated edit script: updates market-price derived columns (H:N)
# across the 8 crafting-job leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to reflect a refreshed Universalis price pull (scheduled runner update).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 193.9  # H12: 197.5 -> 193.9
$ws.Cells.Item(12, 9).Value = 179.5  # I12: 180 -> 179.5
$ws.Cells.Item(12, 10).Value = 197.5  # J12: 200 -> 197.5
$ws.Cells.Item(12, 11).Value = 179.5  # K12: 180 -> 179.5
$ws.Cells.Item(12, 12).Value = 197.5  # L12: 200 -> 197.5
$ws.Cells.Item(12, 13).Value = -9.5  # M12: -10 -> -9.5
$ws.Cells.Item(12, 14).Value = -537.5  # N12: -540 -> -537.5
$ws.Cells.Item(98, 8).Value = 647.4516  # H98: 647.871 -> 647.4516
$ws.Cells.Item(98, 9).Value = 571.4138  # I98: 571.86206 -> 571.4138
$ws.Cells.Item(98, 11).Value = 571.4138  # K98: 571.86206 -> 571.4138
$ws.Cells.Item(98, 13).Value = 926.5862  # M98: 926.13794 -> 926.5862
$ws.Cells.Item(112, 8).Value = 35127.133  # H112: 38885.816 -> 35127.133
$ws.Cells.Item(112, 10).Value = 1671.9286  # J112: 1716.68 -> 1671.9286
$ws.Cells.Item(112, 12).Value = 5015.7858  # L112: 5150.04 -> 5015.7858
$ws.Cells.Item(112, 14).Value = -7231.7858  # N112: -7366.04 -> -7231.7858
$ws.Cells.Item(122, 8).Value = 647.4516  # H122: 647.871 -> 647.4516
$ws.Cells.Item(122, 9).Value = 571.4138  # I122: 571.86206 -> 571.4138
$ws.Cells.Item(122, 11).Value = 1714.2414  # K122: 1715.58618 -> 1714.2414
$ws.Cells.Item(122, 13).Value = 735.7585999999999  # M122: 734.4138199999998 -> 735.7585999999999
$ws.Cells.Item(132, 8).Value = 3234.5405  # H132: 3234.973 -> 3234.5405
$ws.Cells.Item(132, 9).Value = 2824.5151  # I132: 2882.7812 -> 2824.5151
$ws.Cells.Item(132, 10).Value = 6617.25  # J132: 5489 -> 6617.25
$ws.Cells.Item(132, 11).Value = 8473.5453  # K132: 8648.3436 -> 8473.5453
$ws.Cells.Item(132, 12).Value = 19851.75  # L132: 16467 -> 19851.75
$ws.Cells.Item(132, 13).Value = -5943.5453  # M132: -6118.3436 -> -5943.5453
$ws.Cells.Item(132, 14).Value = -24911.75  # N132: -21527 -> -24911.75
$ws.Cells.Item(135, 8).Value = 1066.3214  # H135: 1067.2142 -> 1066.3214
$ws.Cells.Item(135, 9).Value = 448.73914  # I135: 449.82608 -> 448.73914
$ws.Cells.Item(135, 11).Value = 4038.65226  # K135: 4048.43472 -> 4038.65226
$ws.Cells.Item(135, 13).Value = -1503.65226  # M135: -1513.43472 -> -1503.65226

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(22, 8).Value = 1193.4  # H22: 1412.25 -> 1193.4
$ws.Cells.Item(22, 10).Value = 1755.6666  # J22: 2474.5 -> 1755.6666
$ws.Cells.Item(22, 12).Value = 1755.6666  # L22: 2474.5 -> 1755.6666
$ws.Cells.Item(22, 14).Value = -2353.6666  # N22: -3072.5 -> -2353.6666
$ws.Cells.Item(25, 8).Value = 2723  # H25: 5000 -> 2723
$ws.Cells.Item(25, 9).Value = 2723  # I25: 5000 -> 2723
$ws.Cells.Item(25, 11).Value = 2723  # K25: 5000 -> 2723
$ws.Cells.Item(25, 13).Value = -2321  # M25: -4598 -> -2321
$ws.Cells.Item(45, 8).Value = 243792.48  # H45: 280221.66 -> 243792.48
$ws.Cells.Item(45, 9).Value = 371534.28  # I45: 464185 -> 371534.28
$ws.Cells.Item(45, 11).Value = 371534.28  # K45: 464185 -> 371534.28
$ws.Cells.Item(45, 13).Value = -371157.28  # M45: -463808 -> -371157.28
$ws.Cells.Item(61, 8).Value = 1489.8  # H61: 1447.7916 -> 1489.8
$ws.Cells.Item(61, 9).Value = 1260.5  # I61: 1260.591 -> 1260.5
$ws.Cells.Item(61, 10).Value = 3171.3333  # J61: 3507 -> 3171.3333
$ws.Cells.Item(61, 11).Value = 1260.5  # K61: 1260.591 -> 1260.5
$ws.Cells.Item(61, 12).Value = 3171.3333  # L61: 3507 -> 3171.3333
$ws.Cells.Item(61, 13).Value = -1048.5  # M61: -1048.591 -> -1048.5
$ws.Cells.Item(61, 14).Value = -3595.3333  # N61: -3931 -> -3595.3333
$ws.Cells.Item(74, 8).Value = 1579.6774  # H74: 1580.9678 -> 1579.6774
$ws.Cells.Item(74, 9).Value = 1523.4482  # I74: 1524.8276 -> 1523.4482
$ws.Cells.Item(74, 11).Value = 1523.4482  # K74: 1524.8276 -> 1523.4482
$ws.Cells.Item(74, 13).Value = -649.4482  # M74: -650.8276000000001 -> -649.4482
$ws.Cells.Item(77, 8).Value = 1579.6774  # H77: 1580.9678 -> 1579.6774
$ws.Cells.Item(77, 9).Value = 1523.4482  # I77: 1524.8276 -> 1523.4482
$ws.Cells.Item(77, 11).Value = 7617.241  # K77: 7624.138000000001 -> 7617.241
$ws.Cells.Item(77, 13).Value = -3249.241  # M77: -3256.138000000001 -> -3249.241
$ws.Cells.Item(92, 8).Value = 12528887  # H92: 49999 -> 12528887
$ws.Cells.Item(92, 10).Value = 12528887  # J92: 49999 -> 12528887
$ws.Cells.Item(92, 12).Value = 12528887  # L92: 49999 -> 12528887
$ws.Cells.Item(92, 14).Value = -12533879  # N92: -54991 -> -12533879
$ws.Cells.Item(122, 8).Value = 1051.4222  # H122: 1080.711 -> 1051.4222
$ws.Cells.Item(122, 9).Value = 767.4054  # I122: 831.8684 -> 767.4054
$ws.Cells.Item(122, 10).Value = 2365  # J122: 2431.5715 -> 2365
$ws.Cells.Item(122, 11).Value = 2302.2162  # K122: 2495.6052 -> 2302.2162
$ws.Cells.Item(122, 12).Value = 7095  # L122: 7294.7145 -> 7095
$ws.Cells.Item(122, 13).Value = 147.7838000000002  # M122: -45.60519999999997 -> 147.7838000000002
$ws.Cells.Item(122, 14).Value = -11995  # N122: -12194.7145 -> -11995
$ws.Cells.Item(136, 8).Value = 1489.8  # H136: 1447.7916 -> 1489.8
$ws.Cells.Item(136, 9).Value = 1260.5  # I136: 1260.591 -> 1260.5
$ws.Cells.Item(136, 10).Value = 3171.3333  # J136: 3507 -> 3171.3333
$ws.Cells.Item(136, 11).Value = 3781.5  # K136: 3781.773 -> 3781.5
$ws.Cells.Item(136, 12).Value = 9513.999899999999  # L136: 10521 -> 9513.999899999999
$ws.Cells.Item(136, 13).Value = -1231.5  # M136: -1231.773 -> -1231.5
$ws.Cells.Item(136, 14).Value = -14613.9999  # N136: -15621 -> -14613.9999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 10174.409  # H20: 9409.833000000001 -> 10174.409
$ws.Cells.Item(20, 9).Value = 15258.714  # I20: 14308.066 -> 15258.714
$ws.Cells.Item(20, 10).Value = 1276.875  # J20: 1246.1111 -> 1276.875
$ws.Cells.Item(20, 11).Value = 15258.714  # K20: 14308.066 -> 15258.714
$ws.Cells.Item(20, 12).Value = 1276.875  # L20: 1246.1111 -> 1276.875
$ws.Cells.Item(20, 13).Value = -15011.714  # M20: -14061.066 -> -15011.714
$ws.Cells.Item(20, 14).Value = -1770.875  # N20: -1740.1111 -> -1770.875
$ws.Cells.Item(37, 8).Value = 892.2857  # H37: 868.625 -> 892.2857
$ws.Cells.Item(37, 9).Value = 892.2857  # I37: 868.625 -> 892.2857
$ws.Cells.Item(37, 11).Value = 892.2857  # K37: 868.625 -> 892.2857
$ws.Cells.Item(37, 13).Value = -755.2857  # M37: -731.625 -> -755.2857
$ws.Cells.Item(94, 8).Value = 2217.6843  # H94: 2020.6666 -> 2217.6843
$ws.Cells.Item(94, 9).Value = 1650.1818  # I94: 1419.2307 -> 1650.1818
$ws.Cells.Item(94, 11).Value = 1650.1818  # K94: 1419.2307 -> 1650.1818
$ws.Cells.Item(94, 13).Value = -1199.1818  # M94: -968.2307000000001 -> -1199.1818
$ws.Cells.Item(99, 8).Value = 2084.4666  # H99: 2041.625 -> 2084.4666
$ws.Cells.Item(99, 9).Value = 1946.9166  # I99: 1904.7693 -> 1946.9166
$ws.Cells.Item(99, 11).Value = 1946.9166  # K99: 1904.7693 -> 1946.9166
$ws.Cells.Item(99, 13).Value = -448.9166  # M99: -406.7692999999999 -> -448.9166
$ws.Cells.Item(132, 8).Value = 97983.62  # H132: 98213.38 -> 97983.62
$ws.Cells.Item(132, 10).Value = 97983.62  # J132: 98213.38 -> 97983.62
$ws.Cells.Item(132, 12).Value = 97983.62  # L132: 98213.38 -> 97983.62
$ws.Cells.Item(132, 14).Value = -108103.62  # N132: -108333.38 -> -108103.62
$ws.Cells.Item(134, 8).Value = 1580.8889  # H134: 1563.4736 -> 1580.8889
$ws.Cells.Item(134, 9).Value = 1580.8889  # I134: 1563.4736 -> 1580.8889
$ws.Cells.Item(134, 11).Value = 4742.6667  # K134: 4690.4208 -> 4742.6667
$ws.Cells.Item(134, 13).Value = -2207.6667  # M134: -2155.4208 -> -2207.6667

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5901.0356  # H31: 6089.4814 -> 5901.0356
$ws.Cells.Item(31, 9).Value = 5860.3706  # I31: 6054.5 -> 5860.3706
$ws.Cells.Item(31, 11).Value = 5860.3706  # K31: 6054.5 -> 5860.3706
$ws.Cells.Item(31, 13).Value = -5565.3706  # M31: -5759.5 -> -5565.3706
$ws.Cells.Item(34, 8).Value = 5901.0356  # H34: 6089.4814 -> 5901.0356
$ws.Cells.Item(34, 9).Value = 5860.3706  # I34: 6054.5 -> 5860.3706
$ws.Cells.Item(34, 11).Value = 5860.3706  # K34: 6054.5 -> 5860.3706
$ws.Cells.Item(34, 13).Value = -5658.3706  # M34: -5852.5 -> -5658.3706
$ws.Cells.Item(86, 8).Value = 66597  # H86: 66628.78999999999 -> 66597
$ws.Cells.Item(86, 9).Value = 82950.73  # I86: 90690.3 -> 82950.73
$ws.Cells.Item(86, 10).Value = 6633.3335  # J86: 6475 -> 6633.3335
$ws.Cells.Item(86, 11).Value = 82950.73  # K86: 90690.3 -> 82950.73
$ws.Cells.Item(86, 12).Value = 6633.3335  # L86: 6475 -> 6633.3335
$ws.Cells.Item(86, 13).Value = -81827.73  # M86: -89567.3 -> -81827.73
$ws.Cells.Item(86, 14).Value = -8879.333500000001  # N86: -8721 -> -8879.333500000001
$ws.Cells.Item(89, 8).Value = 66597  # H89: 66628.78999999999 -> 66597
$ws.Cells.Item(89, 9).Value = 82950.73  # I89: 90690.3 -> 82950.73
$ws.Cells.Item(89, 10).Value = 6633.3335  # J89: 6475 -> 6633.3335
$ws.Cells.Item(89, 11).Value = 414753.65  # K89: 453451.5 -> 414753.65
$ws.Cells.Item(89, 12).Value = 33166.6675  # L89: 32375 -> 33166.6675
$ws.Cells.Item(89, 13).Value = -409137.65  # M89: -447835.5 -> -409137.65
$ws.Cells.Item(89, 14).Value = -44398.6675  # N89: -43607 -> -44398.6675
$ws.Cells.Item(102, 8).Value = 35500  # H102: 0 -> 35500
$ws.Cells.Item(102, 10).Value = 35500  # J102: 0 -> 35500
$ws.Cells.Item(102, 12).Value = 35500  # L102: 0 -> 35500
$ws.Cells.Item(102, 14).Value = -40368  # N102: None -> -40368
$ws.Cells.Item(103, 8).Value = 12825.818  # H103: 11463.4 -> 12825.818
$ws.Cells.Item(103, 9).Value = 12825.818  # I103: 11463.4 -> 12825.818
$ws.Cells.Item(103, 11).Value = 12825.818  # K103: 11463.4 -> 12825.818
$ws.Cells.Item(103, 13).Value = -11653.818  # M103: -10291.4 -> -11653.818
$ws.Cells.Item(132, 8).Value = 2320.5854  # H132: 2419.975 -> 2320.5854
$ws.Cells.Item(132, 9).Value = 2212.5881  # I132: 2325.182 -> 2212.5881
$ws.Cells.Item(132, 10).Value = 2845.1428  # J132: 2866.8572 -> 2845.1428
$ws.Cells.Item(132, 11).Value = 6637.7643  # K132: 6975.545999999999 -> 6637.7643
$ws.Cells.Item(132, 12).Value = 8535.428400000001  # L132: 8600.571599999999 -> 8535.428400000001
$ws.Cells.Item(132, 13).Value = -4107.7643  # M132: -4445.545999999999 -> -4107.7643
$ws.Cells.Item(132, 14).Value = -13595.4284  # N132: -13660.5716 -> -13595.4284

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 13157965  # H2: 13888962 -> 13157965
$ws.Cells.Item(2, 10).Value = 19230852  # J2: 20833422 -> 19230852
$ws.Cells.Item(2, 12).Value = 115385112  # L2: 125000532 -> 115385112
$ws.Cells.Item(2, 14).Value = -115385338  # N2: -125000758 -> -115385338
$ws.Cells.Item(34, 8).Value = 3999.375  # H34: 4999 -> 3999.375
$ws.Cells.Item(34, 10).Value = 3999.375  # J34: 4999 -> 3999.375
$ws.Cells.Item(34, 12).Value = 11998.125  # L34: 14997 -> 11998.125
$ws.Cells.Item(34, 14).Value = -12166.125  # N34: -15165 -> -12166.125
$ws.Cells.Item(87, 8).Value = 3100  # H87: 3625 -> 3100
$ws.Cells.Item(87, 9).Value = 3100  # I87: 3625 -> 3100
$ws.Cells.Item(87, 11).Value = 9300  # K87: 10875 -> 9300
$ws.Cells.Item(87, 13).Value = -8052  # M87: -9627 -> -8052
$ws.Cells.Item(90, 8).Value = 3100  # H90: 3625 -> 3100
$ws.Cells.Item(90, 9).Value = 3100  # I90: 3625 -> 3100
$ws.Cells.Item(90, 11).Value = 27900  # K90: 32625 -> 27900
$ws.Cells.Item(90, 13).Value = -21660  # M90: -26385 -> -21660
$ws.Cells.Item(104, 8).Value = 50000100  # H104: 25001924 -> 50000100
$ws.Cells.Item(104, 9).Value = 200  # I104: 1349.5 -> 200
$ws.Cells.Item(104, 10).Value = 100000000  # J104: 50002500 -> 100000000
$ws.Cells.Item(104, 11).Value = 600  # K104: 4048.5 -> 600
$ws.Cells.Item(104, 12).Value = 300000000  # L104: 150007500 -> 300000000
$ws.Cells.Item(104, 13).Value = 2021  # M104: -1427.5 -> 2021
$ws.Cells.Item(104, 14).Value = -300005242  # N104: -150012742 -> -300005242
$ws.Cells.Item(131, 8).Value = 16351  # H131: 20328.818 -> 16351
$ws.Cells.Item(131, 9).Value = 849.25  # I131: 819.4 -> 849.25
$ws.Cells.Item(131, 10).Value = 22551.7  # J131: 36586.668 -> 22551.7
$ws.Cells.Item(131, 11).Value = 2547.75  # K131: 2458.2 -> 2547.75
$ws.Cells.Item(131, 12).Value = 67655.10000000001  # L131: 109760.004 -> 67655.10000000001
$ws.Cells.Item(131, 13).Value = 2492.25  # M131: 2581.8 -> 2492.25
$ws.Cells.Item(131, 14).Value = -77735.10000000001  # N131: -119840.004 -> -77735.10000000001
$ws.Cells.Item(140, 8).Value = 2275.5386  # H140: 2331.9167 -> 2275.5386
$ws.Cells.Item(140, 9).Value = 1031  # I140: 1031.1428 -> 1031
$ws.Cells.Item(140, 10).Value = 3727.5  # J140: 4153 -> 3727.5
$ws.Cells.Item(140, 11).Value = 3093  # K140: 3093.4284 -> 3093
$ws.Cells.Item(140, 12).Value = 11182.5  # L140: 12459 -> 11182.5
$ws.Cells.Item(140, 13).Value = 2087  # M140: 2086.5716 -> 2087
$ws.Cells.Item(140, 14).Value = -21542.5  # N140: -22819 -> -21542.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(41, 8).Value = 19999.5  # H41: 20000 -> 19999.5
$ws.Cells.Item(41, 10).Value = 19999.5  # J41: 20000 -> 19999.5
$ws.Cells.Item(41, 12).Value = 19999.5  # L41: 20000 -> 19999.5
$ws.Cells.Item(41, 14).Value = -20709.5  # N41: -20710 -> -20709.5
$ws.Cells.Item(104, 8).Value = 68000  # H104: 34999 -> 68000
$ws.Cells.Item(104, 10).Value = 68000  # J104: 34999 -> 68000
$ws.Cells.Item(104, 12).Value = 68000  # L104: 34999 -> 68000
$ws.Cells.Item(104, 14).Value = -74988  # N104: -41987 -> -74988
$ws.Cells.Item(113, 8).Value = 2128.45  # H113: 2182.6316 -> 2128.45
$ws.Cells.Item(113, 9).Value = 1890  # I113: 1950.8462 -> 1890
$ws.Cells.Item(113, 11).Value = 1890  # K113: 1950.8462 -> 1890
$ws.Cells.Item(113, 13).Value = 280  # M113: 219.1538 -> 280

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(23, 8).Value = 0  # H23: 10006 -> 0
$ws.Cells.Item(23, 9).Value = 0  # I23: 10006 -> 0
$ws.Cells.Item(23, 11).Value = 0  # K23: 10006 -> 0
$ws.Cells.Item(23, 13).ClearContents()  # M23: -9776 -> (empty)
$ws.Cells.Item(40, 8).Value = 2079.8518  # H40: 2083.2593 -> 2079.8518
$ws.Cells.Item(40, 9).Value = 1572  # I40: 1576 -> 1572
$ws.Cells.Item(40, 11).Value = 1572  # K40: 1576 -> 1572
$ws.Cells.Item(40, 13).Value = -1436  # M40: -1440 -> -1436
$ws.Cells.Item(82, 8).Value = 6265.8335  # H82: 4832.625 -> 6265.8335
$ws.Cells.Item(82, 9).Value = 6265.8335  # I82: 4832.625 -> 6265.8335
$ws.Cells.Item(82, 11).Value = 6265.8335  # K82: 4832.625 -> 6265.8335
$ws.Cells.Item(82, 13).Value = -5904.8335  # M82: -4471.625 -> -5904.8335
$ws.Cells.Item(85, 8).Value = 6265.8335  # H85: 4832.625 -> 6265.8335
$ws.Cells.Item(85, 9).Value = 6265.8335  # I85: 4832.625 -> 6265.8335
$ws.Cells.Item(85, 11).Value = 6265.8335  # K85: 4832.625 -> 6265.8335
$ws.Cells.Item(85, 13).Value = -5017.8335  # M85: -3584.625 -> -5017.8335
$ws.Cells.Item(132, 8).Value = 6917.9  # H132: 7097.8423 -> 6917.9
$ws.Cells.Item(132, 9).Value = 3636.2727  # I132: 3650 -> 3636.2727
$ws.Cells.Item(132, 11).Value = 10908.8181  # K132: 10950 -> 10908.8181
$ws.Cells.Item(132, 13).Value = -8378.8181  # M132: -8420 -> -8378.8181
$ws.Cells.Item(136, 8).Value = 2048.5588  # H136: 2104.3428 -> 2048.5588
$ws.Cells.Item(136, 9).Value = 1788.4333  # I136: 1788.4667 -> 1788.4333
$ws.Cells.Item(136, 10).Value = 3999.5  # J136: 3999.6 -> 3999.5
$ws.Cells.Item(136, 11).Value = 5365.2999  # K136: 5365.4001 -> 5365.2999
$ws.Cells.Item(136, 12).Value = 11998.5  # L136: 11998.8 -> 11998.5
$ws.Cells.Item(136, 13).Value = -2815.2999  # M136: -2815.4001 -> -2815.2999
$ws.Cells.Item(136, 14).Value = -17098.5  # N136: -17098.8 -> -17098.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1320.6666  # H122: 1340 -> 1320.6666
$ws.Cells.Item(122, 9).Value = 1280.826  # I122: 1302.6364 -> 1280.826
$ws.Cells.Item(122, 10).Value = 1412.3  # J122: 1422.2 -> 1412.3
$ws.Cells.Item(122, 11).Value = 3842.478  # K122: 3907.9092 -> 3842.478
$ws.Cells.Item(122, 12).Value = 4236.9  # L122: 4266.6 -> 4236.9
$ws.Cells.Item(122, 13).Value = -1392.478  # M122: -1457.9092 -> -1392.478
$ws.Cells.Item(122, 14).Value = -9136.9  # N122: -9166.6 -> -9136.9
$ws.Cells.Item(132, 8).Value = 1696.6724  # H132: 1703.5593 -> 1696.6724
$ws.Cells.Item(132, 9).Value = 1471.5714  # I132: 1482.6492 -> 1471.5714
$ws.Cells.Item(132, 11).Value = 4414.7142  # K132: 4447.9476 -> 4414.7142
$ws.Cells.Item(132, 13).Value = -1884.7142  # M132: -1917.9476 -> -1884.7142
$ws.Cells.Item(136, 8).Value = 1491.4572  # H136: 1482.4166 -> 1491.4572
$ws.Cells.Item(136, 9).Value = 717.6842  # I136: 734.14545 -> 717.6842
$ws.Cells.Item(136, 10).Value = 4884.154  # J136: 3903.2942 -> 4884.154
$ws.Cells.Item(136, 11).Value = 2153.0526  # K136: 2202.43635 -> 2153.0526
$ws.Cells.Item(136, 12).Value = 14652.462  # L136: 11709.8826 -> 14652.462
$ws.Cells.Item(136, 13).Value = 396.9474  # M136: 347.5636500000001 -> 396.9474
$ws.Cells.Item(136, 14).Value = -19752.462  # N136: -16809.8826 -> -19752.462
$ws.Cells.Item(141, 8).Value = 0  # H141: 94990 -> 0
$ws.Cells.Item(141, 10).Value = 0  # J141: 94990 -> 0
$ws.Cells.Item(141, 12).Value = 0  # L141: 94990 -> 0
$ws.Cells.Item(141, 14).ClearContents()  # N141: -105350 -> (empty)
